$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 210.6
$ws.Range("I33").Value = 210.6
$ws.Range("K33").Value = 210.6
$ws.Range("M33").Value = 18.40000000000001

# Row 53
$ws.Range("H53").Value = 334.6316
$ws.Range("J53").Value = 398.7
$ws.Range("L53").Value = 398.7
$ws.Range("N53").Value = -1672.7

# Row 74
$ws.Range("H74").Value = 27837170
$ws.Range("I74").Value = 27837170
$ws.Range("K74").Value = 27837170
$ws.Range("M74").Value = -27836234

# Row 77
$ws.Range("H77").Value = 27837170
$ws.Range("I77").Value = 27837170
$ws.Range("K77").Value = 139185850
$ws.Range("M77").Value = -139181170

# Row 92
$ws.Range("H92").Value = 902
$ws.Range("I92").Value = 791.75
$ws.Range("J92").Value = 1067.375
$ws.Range("K92").Value = 791.75
$ws.Range("L92").Value = 1067.375
$ws.Range("M92").Value = 456.25
$ws.Range("N92").Value = -3563.375

# Row 96
$ws.Range("H96").Value = 1309.2
$ws.Range("I96").Value = 930
$ws.Range("K96").Value = 2790
$ws.Range("M96").Value = -1417

# Row 100
$ws.Range("H100").Value = 2883.1667
$ws.Range("I100").Value = 2979.8
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 2979.8
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -2438.8
$ws.Range("N100").Value = -3482

# Row 112
$ws.Range("H112").Value = 49548.188
$ws.Range("I112").Value = 85271
$ws.Range("J112").Value = 35720
$ws.Range("K112").Value = 255813
$ws.Range("L112").Value = 107160
$ws.Range("M112").Value = -254705
$ws.Range("N112").Value = -109376

# Row 116
$ws.Range("H116").Value = 5333.524
$ws.Range("I116").Value = 5036
$ws.Range("K116").Value = 5036
$ws.Range("M116").Value = -1594

# Row 135
$ws.Range("H135").Value = 115385010
$ws.Range("I135").Value = 50000450
$ws.Range("K135").Value = 450004050
$ws.Range("M135").Value = -450001515

# Row 138
$ws.Range("H138").Value = 3106.926
$ws.Range("J138").Value = 3799.7334
$ws.Range("L138").Value = 11399.2002
$ws.Range("N138").Value = -21679.2002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2868
$ws.Range("I32").Value = 2868
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2868
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2581
$ws.Range("N32").Value = $null

# Row 45
$ws.Range("H45").Value = 1605.75
$ws.Range("I45").Value = 1341.1666
$ws.Range("K45").Value = 1341.1666
$ws.Range("M45").Value = -964.1666

# Row 97
$ws.Range("H97").Value = 654.11536
$ws.Range("I97").Value = 626.4545000000001
$ws.Range("J97").Value = 806.25
$ws.Range("K97").Value = 626.4545000000001
$ws.Range("L97").Value = 806.25
$ws.Range("M97").Value = -130.4545000000001
$ws.Range("N97").Value = -1798.25

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 2231.4827
$ws.Range("I122").Value = 2218.6
$ws.Range("J122").Value = 2312
$ws.Range("K122").Value = 6655.799999999999
$ws.Range("L122").Value = 6936
$ws.Range("M122").Value = -4205.799999999999
$ws.Range("N122").Value = -11836

# Row 134
$ws.Range("H134").Value = 16668574
$ws.Range("I134").Value = 19232584
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 57697752
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -57695217
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 82576.664
$ws.Range("J121").Value = 21826.8
$ws.Range("L121").Value = 65480.39999999999
$ws.Range("N121").Value = -68100.39999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 679.4666999999999
$ws.Range("I97").Value = 418
$ws.Range("J97").Value = 1131.091
$ws.Range("K97").Value = 418
$ws.Range("L97").Value = 1131.091
$ws.Range("M97").Value = 78
$ws.Range("N97").Value = -2123.091

# Row 102
$ws.Range("H102").Value = 1199.25
$ws.Range("I102").Value = 888.44446
$ws.Range("J102").Value = 2131.6667
$ws.Range("K102").Value = 888.44446
$ws.Range("L102").Value = 2131.6667
$ws.Range("M102").Value = 733.55554
$ws.Range("N102").Value = -5375.6667

# Row 113
$ws.Range("H113").Value = 79491.84
$ws.Range("I113").Value = 101839.7
$ws.Range("K113").Value = 101839.7
$ws.Range("M113").Value = -99669.7

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2373.5833
$ws.Range("I7").Value = 1848.3
$ws.Range("K7").Value = 1848.3
$ws.Range("M7").Value = -1736.3

# Row 22
$ws.Range("H22").Value = 3943.75
$ws.Range("I22").Value = 3935.7144
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 3935.7144
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -3640.7144
$ws.Range("N22").Value = -4590

# Row 27
$ws.Range("H27").Value = 3943.75
$ws.Range("I27").Value = 3935.7144
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 3935.7144
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -3828.7144
$ws.Range("N27").Value = -4214

# Row 40
$ws.Range("H40").Value = 4797.385
$ws.Range("I40").Value = 4549.6
$ws.Range("J40").Value = 5623.3335
$ws.Range("K40").Value = 4549.6
$ws.Range("L40").Value = 5623.3335
$ws.Range("M40").Value = -4413.6
$ws.Range("N40").Value = -5895.3335

# Row 55
$ws.Range("H55").Value = 236.33333
$ws.Range("J55").Value = 259.4
$ws.Range("L55").Value = 259.4
$ws.Range("N55").Value = -605.4

# Row 122
$ws.Range("H122").Value = 5310.294
$ws.Range("I122").Value = 5310.294
$ws.Range("K122").Value = 15930.882
$ws.Range("M122").Value = -13480.882

# Row 126
$ws.Range("H126").Value = 2373.5833
$ws.Range("I126").Value = 1848.3
$ws.Range("K126").Value = 5544.9
$ws.Range("M126").Value = -3074.9

# Row 132
$ws.Range("H132").Value = 27791744
$ws.Range("I132").Value = 27791744
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 83375232
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -83372702
$ws.Range("N132").Value = $null

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

# Row 136
$ws.Range("H136").Value = 2791.5557
$ws.Range("I136").Value = 2066.3333
$ws.Range("J136").Value = 2998.762
$ws.Range("K136").Value = 6198.999899999999
$ws.Range("L136").Value = 8996.286
$ws.Range("M136").Value = -3648.999899999999
$ws.Range("N136").Value = -14096.286

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null

# Row 81
$ws.Range("H81").Value = 3425.8333
$ws.Range("I81").Value = 3425.8333
$ws.Range("K81").Value = 6851.6666
$ws.Range("M81").Value = -5790.6666

# Row 84
$ws.Range("H84").Value = 3425.8333
$ws.Range("I84").Value = 3425.8333
$ws.Range("K84").Value = 34258.333
$ws.Range("M84").Value = -28954.333

# Row 122
$ws.Range("H122").Value = 5409.607
$ws.Range("I122").Value = 5219.24
$ws.Range("K122").Value = 15657.72
$ws.Range("M122").Value = -13207.72

# Row 126
$ws.Range("H126").Value = 1322.6333
$ws.Range("I126").Value = 1287.5555
$ws.Range("J126").Value = 1638.3334
$ws.Range("K126").Value = 3862.6665
$ws.Range("L126").Value = 4915.0002
$ws.Range("M126").Value = -1392.6665
$ws.Range("N126").Value = -9855.0002
